# Insert a new weekly price record as row 305 ("Fruta / hortaliza, semanal").
# This shifts the existing rows 305:387 down to 306:388 (Excel preserves the
# formatting of the row being pushed down, including the date number format
# on column D), and grows the sheet from A1:R387 to A1:R388.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(305).Insert()

$ws.Range("A305").Value = 8
$ws.Range("B305").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C305").Value = 'Coquimbo'
$ws.Range("D305").Value = 44988
$ws.Range("E305").Value = 4
$ws.Range("F305").Value = 100112021
$ws.Range("G305").Value = 'Ají'
$ws.Range("H305").Value = 'Inferno'
$ws.Range("I305").Value = 'Primera'
$ws.Range("J305").Value = 480
$ws.Range("K305").Value = 10000
$ws.Range("L305").Value = 11000
$ws.Range("M305").Value = 10500
$ws.Range("N305").Value = '$/caja 15 kilos'
$ws.Range("O305").Value = 'Provincia de Limarí'
$ws.Range("P305").Value = 700
$ws.Range("Q305").Value = 15
$ws.Range("R305").Value = 'Hortaliza'
